$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Enhed"
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "By"

$ws.Range("B11").Select()
